$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "27.145.82"
$ws.Cells.Item(2, 5).Value = "  -0.04%  "

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.900.16"
$ws.Cells.Item(3, 5).Value = "  -0.04%  "

$ws.Cells.Item(4, 5).Value = "  +0.10%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "307.09"
$ws.Cells.Item(5, 5).Value = "  +0.18%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "1.003"
$ws.Cells.Item(6, 5).Value = "  +0.23%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.5235"
$ws.Cells.Item(7, 5).Value = "  -0.09%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3806"
$ws.Cells.Item(8, 5).Value = "  +0.93%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.07287"

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "21.30"
$ws.Cells.Item(10, 5).Value = "  +0.74%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.9042"
$ws.Cells.Item(11, 5).Value = "  +0.57%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.08198"
$ws.Cells.Item(12, 5).Value = "  -2.18%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "1.892.61"
$ws.Cells.Item(13, 5).Value = "  -0.48%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "95.47"
$ws.Cells.Item(14, 5).Value = "  +0.86%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "5.354"
$ws.Cells.Item(15, 5).Value = "  +1.66%  "

$ws.Cells.Item(16, 5).Value = "  +0.10%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.000008658"
$ws.Cells.Item(17, 5).Value = "  +0.86%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "14.67"
$ws.Cells.Item(18, 5).Value = "  +1.21%  "

$ws.Cells.Item(19, 5).Value = "  +0.18%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "27.190.94"
$ws.Cells.Item(20, 5).Value = "  -0.03%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "5.123"
$ws.Cells.Item(21, 5).Value = "  +1.23%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "2.127.04"
$ws.Cells.Item(22, 5).Value = "  -0.81%  "

$ws.Cells.Item(23, 5).Value = "  +1.94%  "

$ws.Cells.Item(24, 5).Value = "  +0.78%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "149.37"
$ws.Cells.Item(25, 5).Value = "  +1.96%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.315"
$ws.Cells.Item(26, 5).Value = "  +1.50%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "18.28"
$ws.Cells.Item(27, 5).Value = "  +0.97%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "1.743"
$ws.Cells.Item(28, 5).Value = "  -0.78%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "115.40"
$ws.Cells.Item(29, 5).Value = "  +0.63%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "4.821"
$ws.Cells.Item(30, 5).Value = "  +0.83%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "4.890"
$ws.Cells.Item(31, 5).Value = "  -0.64%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.09219"
$ws.Cells.Item(32, 5).Value = "  -0.11%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.05040"
$ws.Cells.Item(33, 5).Value = "  -0.17%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.7926"
$ws.Cells.Item(34, 5).Value = "  -2.79%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.222"
$ws.Cells.Item(35, 5).Value = "  -1.07%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "2.979"
$ws.Cells.Item(36, 5).Value = "  +0.66%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "3.345"
$ws.Cells.Item(37, 5).Value = "  -0.76%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "2.636"
$ws.Cells.Item(38, 5).Value = "  +2.61%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.5732"
$ws.Cells.Item(39, 5).Value = "  +0.73%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.01992"
$ws.Cells.Item(40, 5).Value = "  +1.06%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "1.082"
$ws.Cells.Item(41, 5).Value = "  +0.73%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "9.035"
$ws.Cells.Item(42, 5).Value = "  +1.07%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "6.612"
$ws.Cells.Item(43, 5).Value = "  -0.69%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "116.50"
$ws.Cells.Item(44, 5).Value = "  -1.62%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.1516"
$ws.Cells.Item(45, 5).Value = "  +0.42%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.4897"
$ws.Cells.Item(46, 5).Value = "  +1.61%  "

$ws.Cells.Item(47, 2).Value = "PaxDollar"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "1.004"
$ws.Cells.Item(47, 5).Value = "  +0.28%  "

$ws.Cells.Item(48, 2).Value = "EnergySwap"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "10.16"
$ws.Cells.Item(48, 5).Value = "  -0.14%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "1.638"
$ws.Cells.Item(49, 5).Value = "  +1.84%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "38.55"
$ws.Cells.Item(50, 5).Value = "  +3.02%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "64.11"
$ws.Cells.Item(51, 5).Value = "  +0.92%  "
